# Apply "Started" column (C) Yes/No toggles on the "by Coach" sheet,
# and update the saved view (frozen pane top-left cell + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Rows in column C whose value needs to be flipped between "Yes" and "No".
$rows = @(3,5,8,10,11,12,40,46,51,54,55,60,63,64,66,67,69,76,77,83,84)

foreach ($r in $rows) {
    $cell = $ws.Range("C$r")
    if ($cell.Value2 -eq "Yes") {
        $cell.Value = "No"
    } else {
        $cell.Value = "Yes"
    }
}

# Update the frozen-pane top-left cell and the active selection to match
# the new scroll position / selected cell recorded in the saved view.
# Re-freeze the pane anchored at the new top-left cell (row 63) so the
# view's top-left cell matches the saved file, then move the selection
# to the new active cell.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A63").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("C85").Select()
